$wb = $excel.ActiveWorkbook

# Sheet: ALC (165 cell updates)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 257.6111
$ws.Range("J2").Value = 303.8
$ws.Range("L2").Value = 303.8
$ws.Range("N2").Value = -529.8
$ws.Range("H4").Value = 196.85715
$ws.Range("I4").Value = 213
$ws.Range("K4").Value = 213
$ws.Range("M4").Value = -99
$ws.Range("H5").Value = 250000140
$ws.Range("I5").Value = 194.66667
$ws.Range("K5").Value = 194.66667
$ws.Range("M5").Value = -79.66667000000001
$ws.Range("H6").Value = 10015
$ws.Range("I6").Value = 10015
$ws.Range("K6").Value = 30045
$ws.Range("M6").Value = -29933
$ws.Range("H9").Value = 271.1613
$ws.Range("I9").Value = 78.13636
$ws.Range("J9").Value = 743
$ws.Range("K9").Value = 78.13636
$ws.Range("L9").Value = 743
$ws.Range("M9").Value = 90.86364
$ws.Range("N9").Value = -1081
$ws.Range("H12").Value = 134
$ws.Range("I12").Value = 134
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 134
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = 36
$ws.Range("N12").ClearContents()
$ws.Range("H15").Value = 2266.1177
$ws.Range("I15").Value = 2266.1177
$ws.Range("K15").Value = 6798.353099999999
$ws.Range("M15").Value = -6629.353099999999
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("M21").ClearContents()
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("M23").ClearContents()
$ws.Range("H38").Value = 213.44444
$ws.Range("I38").Value = 213.44444
$ws.Range("K38").Value = 640.33332
$ws.Range("M38").Value = -268.33332
$ws.Range("H48").Value = 999
$ws.Range("I48").Value = 999
$ws.Range("J48").Value = 0
$ws.Range("K48").Value = 2997
$ws.Range("L48").Value = 0
$ws.Range("M48").Value = -2705
$ws.Range("N48").ClearContents()
$ws.Range("H55").Value = 185.58333
$ws.Range("J55").Value = 400.5
$ws.Range("L55").Value = 400.5
$ws.Range("N55").Value = -828.5
$ws.Range("H56").Value = 999
$ws.Range("I56").Value = 999
$ws.Range("J56").Value = 0
$ws.Range("K56").Value = 2997
$ws.Range("L56").Value = 0
$ws.Range("M56").Value = -2463
$ws.Range("N56").ClearContents()
$ws.Range("H58").Value = 662.5
$ws.Range("I58").Value = 662.5
$ws.Range("K58").Value = 1987.5
$ws.Range("M58").Value = -1837.5
$ws.Range("H74").Value = 14229
$ws.Range("I74").Value = 9238.5
$ws.Range("K74").Value = 9238.5
$ws.Range("M74").Value = -8302.5
$ws.Range("H76").Value = 15418.9
$ws.Range("I76").Value = 26339
$ws.Range("J76").Value = 4498.8
$ws.Range("K76").Value = 26339
$ws.Range("L76").Value = 4498.8
$ws.Range("M76").Value = -26024
$ws.Range("N76").Value = -5128.8
$ws.Range("H77").Value = 14229
$ws.Range("I77").Value = 9238.5
$ws.Range("K77").Value = 46192.5
$ws.Range("M77").Value = -41512.5
$ws.Range("H79").Value = 15418.9
$ws.Range("I79").Value = 26339
$ws.Range("J79").Value = 4498.8
$ws.Range("K79").Value = 26339
$ws.Range("L79").Value = 4498.8
$ws.Range("M79").Value = -25247
$ws.Range("N79").Value = -6682.8
$ws.Range("H80").Value = 998.9318
$ws.Range("I80").Value = 1040.4722
$ws.Range("J80").Value = 812
$ws.Range("K80").Value = 3121.4166
$ws.Range("L80").Value = 2436
$ws.Range("M80").Value = -2123.4166
$ws.Range("N80").Value = -4432
$ws.Range("H83").Value = 998.9318
$ws.Range("I83").Value = 1040.4722
$ws.Range("J83").Value = 812
$ws.Range("K83").Value = 9364.2498
$ws.Range("L83").Value = 7308
$ws.Range("M83").Value = -4372.2498
$ws.Range("N83").Value = -17292
$ws.Range("H86").Value = 2354.5557
$ws.Range("I86").Value = 2897.8
$ws.Range("K86").Value = 2897.8
$ws.Range("M86").Value = -1774.8
$ws.Range("H89").Value = 2354.5557
$ws.Range("I89").Value = 2897.8
$ws.Range("K89").Value = 14489
$ws.Range("M89").Value = -8873
$ws.Range("H97").Value = 2400.875
$ws.Range("J97").Value = 2400.875
$ws.Range("L97").Value = 7202.625
$ws.Range("N97").Value = -8194.625
$ws.Range("H106").Value = 1969.3334
$ws.Range("I106").Value = 1969.3334
$ws.Range("K106").Value = 1969.3334
$ws.Range("M106").Value = -1338.3334
$ws.Range("H107").Value = 26317982
$ws.Range("I107").Value = 16668071
$ws.Range("J107").Value = 62505148
$ws.Range("K107").Value = 16668071
$ws.Range("L107").Value = 62505148
$ws.Range("M107").Value = -16666151
$ws.Range("N107").Value = -62508988
$ws.Range("H108").Value = 78333
$ws.Range("J108").Value = 87499.5
$ws.Range("L108").Value = 87499.5
$ws.Range("N108").Value = -95179.5
$ws.Range("H112").Value = 2021763.8
$ws.Range("J112").Value = 2211160.5
$ws.Range("L112").Value = 6633481.5
$ws.Range("N112").Value = -6635697.5
$ws.Range("H129").Value = 3222.2307
$ws.Range("J129").Value = 3911
$ws.Range("L129").Value = 11733
$ws.Range("N129").Value = -21733
$ws.Range("H132").Value = 2602.4482
$ws.Range("I132").Value = 2655.4783
$ws.Range("J132").Value = 2399.1667
$ws.Range("K132").Value = 7966.4349
$ws.Range("L132").Value = 7197.500100000001
$ws.Range("M132").Value = -5436.4349
$ws.Range("N132").Value = -12257.5001
$ws.Range("H133").Value = 126666.336
$ws.Range("J133").Value = 126666.336
$ws.Range("L133").Value = 126666.336
$ws.Range("N133").Value = -136786.336
$ws.Range("H137").Value = 6669.375
$ws.Range("I137").Value = 8785.066000000001
$ws.Range("K137").Value = 26355.198
$ws.Range("M137").Value = -23805.198
$ws.Range("H138").Value = 5892.16
$ws.Range("I138").Value = 3942
$ws.Range("J138").Value = 5973.4165
$ws.Range("K138").Value = 11826
$ws.Range("L138").Value = 17920.2495
$ws.Range("M138").Value = -6686
$ws.Range("N138").Value = -28200.2495
$ws.Range("H141").Value = 2742.0286
$ws.Range("I141").Value = 2654.8276
$ws.Range("K141").Value = 7964.4828
$ws.Range("M141").Value = -2784.4828

# Sheet: ARM (69 cell updates)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 689
$ws.Range("I4").Value = 689
$ws.Range("K4").Value = 689
$ws.Range("M4").Value = -573
$ws.Range("H21").Value = 2216
$ws.Range("I21").Value = 114
$ws.Range("J21").Value = 3267
$ws.Range("K21").Value = 114
$ws.Range("L21").Value = 3267
$ws.Range("M21").Value = 260
$ws.Range("N21").Value = -4015
$ws.Range("H28").Value = 15850.5
$ws.Range("I28").Value = 5302.3335
$ws.Range("J28").Value = 47495
$ws.Range("K28").Value = 5302.3335
$ws.Range("L28").Value = 47495
$ws.Range("M28").Value = -5110.3335
$ws.Range("N28").Value = -47879
$ws.Range("H32").Value = 21995.426
$ws.Range("I32").Value = 18805.45
$ws.Range("K32").Value = 18805.45
$ws.Range("M32").Value = -18518.45
$ws.Range("H61").Value = 6982.8335
$ws.Range("I61").Value = 5366.8887
$ws.Range("K61").Value = 5366.8887
$ws.Range("M61").Value = -5154.8887
$ws.Range("H88").Value = 2009.4166
$ws.Range("I88").Value = 1971.125
$ws.Range("J88").Value = 2028.5625
$ws.Range("K88").Value = 1971.125
$ws.Range("L88").Value = 2028.5625
$ws.Range("M88").Value = -1565.125
$ws.Range("N88").Value = -2840.5625
$ws.Range("H91").Value = 2009.4166
$ws.Range("I91").Value = 1971.125
$ws.Range("J91").Value = 2028.5625
$ws.Range("K91").Value = 1971.125
$ws.Range("L91").Value = 2028.5625
$ws.Range("M91").Value = -567.125
$ws.Range("N91").Value = -4836.5625
$ws.Range("H99").Value = 15850.5
$ws.Range("I99").Value = 5302.3335
$ws.Range("J99").Value = 47495
$ws.Range("K99").Value = 5302.3335
$ws.Range("L99").Value = 47495
$ws.Range("M99").Value = -2307.3335
$ws.Range("N99").Value = -53485
$ws.Range("H122").Value = 7601.4546
$ws.Range("I122").Value = 4087
$ws.Range("J122").Value = 11818.8
$ws.Range("K122").Value = 12261
$ws.Range("L122").Value = 35456.39999999999
$ws.Range("M122").Value = -9811
$ws.Range("N122").Value = -40356.39999999999
$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").ClearContents()
$ws.Range("H132").Value = 5256.125
$ws.Range("I132").Value = 5618.6665
$ws.Range("J132").Value = 4893.5835
$ws.Range("K132").Value = 16855.9995
$ws.Range("L132").Value = 14680.7505
$ws.Range("M132").Value = -14325.9995
$ws.Range("N132").Value = -19740.7505
$ws.Range("H136").Value = 6982.8335
$ws.Range("I136").Value = 5366.8887
$ws.Range("K136").Value = 16100.6661
$ws.Range("M136").Value = -13550.6661

# Sheet: BSM (19 cell updates)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2584.6943
$ws.Range("I20").Value = 1505.9131
$ws.Range("K20").Value = 1505.9131
$ws.Range("M20").Value = -1258.9131
$ws.Range("H107").Value = 2167.3809
$ws.Range("I107").Value = 1945.3125
$ws.Range("K107").Value = 1945.3125
$ws.Range("M107").Value = -25.3125
$ws.Range("H134").Value = 3878.946
$ws.Range("I134").Value = 5390.3335
$ws.Range("J134").Value = 2447.1052
$ws.Range("K134").Value = 16171.0005
$ws.Range("L134").Value = 7341.3156
$ws.Range("M134").Value = -13636.0005
$ws.Range("N134").Value = -12411.3156
$ws.Range("H141").Value = 58749.75
$ws.Range("J141").Value = 58749.75
$ws.Range("L141").Value = 58749.75
$ws.Range("N141").Value = -69109.75

# Sheet: CRP (66 cell updates)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 3750
$ws.Range("I6").Value = 5500
$ws.Range("K6").Value = 5500
$ws.Range("M6").Value = -5387
$ws.Range("H7").Value = 96.75
$ws.Range("I7").Value = 103.42857
$ws.Range("J7").Value = 50
$ws.Range("K7").Value = 103.42857
$ws.Range("L7").Value = 50
$ws.Range("M7").Value = 9.571430000000007
$ws.Range("N7").Value = -276
$ws.Range("H25").Value = 3750
$ws.Range("J25").Value = 3750
$ws.Range("L25").Value = 3750
$ws.Range("N25").Value = -4098
$ws.Range("H31").Value = 35383424
$ws.Range("I31").Value = 9807282
$ws.Range("J31").Value = 62558076
$ws.Range("K31").Value = 9807282
$ws.Range("L31").Value = 62558076
$ws.Range("M31").Value = -9806987
$ws.Range("N31").Value = -62558666
$ws.Range("H32").Value = 20301.8
$ws.Range("J32").Value = 24999
$ws.Range("L32").Value = 24999
$ws.Range("N32").Value = -25631
$ws.Range("H34").Value = 35383424
$ws.Range("I34").Value = 9807282
$ws.Range("J34").Value = 62558076
$ws.Range("K34").Value = 9807282
$ws.Range("L34").Value = 62558076
$ws.Range("M34").Value = -9807080
$ws.Range("N34").Value = -62558480
$ws.Range("H58").Value = 3568.5
$ws.Range("I58").Value = 3228.2
$ws.Range("J58").Value = 4135.6665
$ws.Range("K58").Value = 3228.2
$ws.Range("L58").Value = 4135.6665
$ws.Range("M58").Value = -3025.2
$ws.Range("N58").Value = -4541.6665
$ws.Range("H94").Value = 1264.48
$ws.Range("J94").Value = 1288.5
$ws.Range("L94").Value = 1288.5
$ws.Range("N94").Value = -2190.5
$ws.Range("H122").Value = 33282
$ws.Range("I122").Value = 2025.75
$ws.Range("K122").Value = 6077.25
$ws.Range("M122").Value = -3627.25
$ws.Range("H132").Value = 3467.0698
$ws.Range("I132").Value = 2302.742
$ws.Range("K132").Value = 6908.226000000001
$ws.Range("M132").Value = -4378.226000000001
$ws.Range("H134").Value = 5416.697
$ws.Range("I134").Value = 6225.923
$ws.Range("J134").Value = 2411
$ws.Range("K134").Value = 18677.769
$ws.Range("L134").Value = 7233
$ws.Range("M134").Value = -16142.769
$ws.Range("N134").Value = -12303
$ws.Range("H136").Value = 3568.5
$ws.Range("I136").Value = 3228.2
$ws.Range("J136").Value = 4135.6665
$ws.Range("K136").Value = 9684.599999999999
$ws.Range("L136").Value = 12406.9995
$ws.Range("M136").Value = -7134.599999999999
$ws.Range("N136").Value = -17506.9995

# Sheet: CUL (95 cell updates)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 2266241
$ws.Range("I4").Value = 2302419.5
$ws.Range("J4").Value = 1000000
$ws.Range("K4").Value = 6907258.5
$ws.Range("L4").Value = 3000000
$ws.Range("M4").Value = -6907146.5
$ws.Range("N4").Value = -3000224
$ws.Range("H5").Value = 948.5
$ws.Range("I5").Value = 1002.9091
$ws.Range("J5").Value = 350
$ws.Range("K5").Value = 3008.7273
$ws.Range("L5").Value = 1050
$ws.Range("M5").Value = -2896.7273
$ws.Range("N5").Value = -1274
$ws.Range("H7").Value = 126532.5
$ws.Range("I7").Value = 200354
$ws.Range("K7").Value = 601062
$ws.Range("M7").Value = -600950
$ws.Range("H10").Value = 1884.2222
$ws.Range("I10").Value = 1035.6
$ws.Range("K10").Value = 3106.8
$ws.Range("M10").Value = -2967.8
$ws.Range("H37").Value = 178467.6
$ws.Range("J37").Value = 178467.6
$ws.Range("L37").Value = 535402.8
$ws.Range("N37").Value = -535626.8
$ws.Range("H40").Value = 93.375
$ws.Range("I40").Value = 84.5
$ws.Range("K40").Value = 338
$ws.Range("M40").Value = -269
$ws.Range("H47").Value = 12482.615
$ws.Range("I47").Value = 2059.8
$ws.Range("J47").Value = 18996.875
$ws.Range("K47").Value = 6179.400000000001
$ws.Range("L47").Value = 56990.625
$ws.Range("M47").Value = -5748.400000000001
$ws.Range("N47").Value = -57852.625
$ws.Range("H92").Value = 784.6667
$ws.Range("I92").Value = 848.5
$ws.Range("J92").Value = 752.75
$ws.Range("K92").Value = 2545.5
$ws.Range("L92").Value = 2258.25
$ws.Range("M92").Value = -1297.5
$ws.Range("N92").Value = -4754.25
$ws.Range("H97").Value = 403.57144
$ws.Range("I97").Value = 216.5
$ws.Range("J97").Value = 478.4
$ws.Range("K97").Value = 649.5
$ws.Range("L97").Value = 1435.2
$ws.Range("M97").Value = -153.5
$ws.Range("N97").Value = -2427.2
$ws.Range("H107").Value = 684.9756
$ws.Range("I107").Value = 1041.3684
$ws.Range("J107").Value = 577.49207
$ws.Range("K107").Value = 3124.1052
$ws.Range("L107").Value = 1732.47621
$ws.Range("M107").Value = -1204.1052
$ws.Range("N107").Value = -5572.47621
$ws.Range("H129").Value = 74654600
$ws.Range("I129").Value = 135417340
$ws.Range("K129").Value = 406252020
$ws.Range("M129").Value = -406247020
$ws.Range("H131").Value = 16950554
$ws.Range("J131").Value = 1465.6604
$ws.Range("L131").Value = 4396.9812
$ws.Range("N131").Value = -14476.9812
$ws.Range("H134").Value = 5608.5557
$ws.Range("I134").Value = 5934.625
$ws.Range("J134").Value = 3000
$ws.Range("K134").Value = 17803.875
$ws.Range("L134").Value = 9000
$ws.Range("M134").Value = -12733.875
$ws.Range("N134").Value = -19140
$ws.Range("H135").Value = 948.5
$ws.Range("I135").Value = 1002.9091
$ws.Range("J135").Value = 350
$ws.Range("K135").Value = 9026.1819
$ws.Range("L135").Value = 3150
$ws.Range("M135").Value = -6491.1819
$ws.Range("N135").Value = -8220
$ws.Range("H136").Value = 1575.5238
$ws.Range("I136").Value = 1523.8536
$ws.Range("K136").Value = 4571.560799999999
$ws.Range("M136").Value = 528.4392000000007
$ws.Range("H137").Value = 60610144
$ws.Range("I137").Value = 5799.2
$ws.Range("J137").Value = 111113770
$ws.Range("K137").Value = 17397.6
$ws.Range("L137").Value = 333341310
$ws.Range("M137").Value = -12297.6
$ws.Range("N137").Value = -333351510
$ws.Range("H139").Value = 2125.1
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

# Sheet: GSM (45 cell updates)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 108635.21
$ws.Range("I2").Value = 325047.88
$ws.Range("K2").Value = 325047.88
$ws.Range("M2").Value = -324934.88
$ws.Range("H12").Value = 4900.8
$ws.Range("I12").Value = 4912
$ws.Range("J12").Value = 4884
$ws.Range("K12").Value = 4912
$ws.Range("L12").Value = 4884
$ws.Range("M12").Value = -4772
$ws.Range("N12").Value = -5164
$ws.Range("H23").Value = 2991
$ws.Range("I23").Value = 2991
$ws.Range("K23").Value = 2991
$ws.Range("M23").Value = -2768
$ws.Range("H36").Value = 5244.5
$ws.Range("I36").Value = 4661
$ws.Range("K36").Value = 4661
$ws.Range("M36").Value = -4176
$ws.Range("H70").Value = 58828852
$ws.Range("I70").Value = 4198.8
$ws.Range("J70").Value = 142864060
$ws.Range("K70").Value = 4198.8
$ws.Range("L70").Value = 142864060
$ws.Range("M70").Value = -3928.8
$ws.Range("N70").Value = -142864600
$ws.Range("H73").Value = 58828852
$ws.Range("I73").Value = 4198.8
$ws.Range("J73").Value = 142864060
$ws.Range("K73").Value = 4198.8
$ws.Range("L73").Value = 142864060
$ws.Range("M73").Value = -3262.8
$ws.Range("N73").Value = -142865932
$ws.Range("H127").Value = 44999
$ws.Range("J127").Value = 44999
$ws.Range("L127").Value = 44999
$ws.Range("N127").Value = -54919
$ws.Range("H132").Value = 39317.484
$ws.Range("I132").Value = 67800.56
$ws.Range("K132").Value = 203401.68
$ws.Range("M132").Value = -200871.68
$ws.Range("H141").Value = 92291.664
$ws.Range("I141").Value = 45454.547
$ws.Range("K141").Value = 45454.547
$ws.Range("M141").Value = -40274.547

# Sheet: LTW (62 cell updates)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H10").Value = 999
$ws.Range("J10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("N10").ClearContents()
$ws.Range("H16").Value = 41680252
$ws.Range("I16").Value = 55559332
$ws.Range("K16").Value = 55559332
$ws.Range("M16").Value = -55559162
$ws.Range("H22").Value = 2860.2856
$ws.Range("I22").Value = 2222.8572
$ws.Range("J22").Value = 3497.7144
$ws.Range("K22").Value = 2222.8572
$ws.Range("L22").Value = 3497.7144
$ws.Range("M22").Value = -1927.8572
$ws.Range("N22").Value = -4087.7144
$ws.Range("H27").Value = 2860.2856
$ws.Range("I27").Value = 2222.8572
$ws.Range("J27").Value = 3497.7144
$ws.Range("K27").Value = 2222.8572
$ws.Range("L27").Value = 3497.7144
$ws.Range("M27").Value = -2115.8572
$ws.Range("N27").Value = -3711.7144
$ws.Range("H33").Value = 27874.5
$ws.Range("I33").Value = 29749
$ws.Range("K33").Value = 29749
$ws.Range("M33").Value = -29459
$ws.Range("H34").Value = 7299
$ws.Range("I34").Value = 7873.75
$ws.Range("J34").Value = 5000
$ws.Range("K34").Value = 7873.75
$ws.Range("L34").Value = 5000
$ws.Range("M34").Value = -7701.75
$ws.Range("N34").Value = -5344
$ws.Range("H40").Value = 687820.2
$ws.Range("I40").Value = 1023939.94
$ws.Range("K40").Value = 1023939.94
$ws.Range("M40").Value = -1023803.94
$ws.Range("H99").Value = 34666.332
$ws.Range("I99").Value = 27000
$ws.Range("J99").Value = 49999
$ws.Range("K99").Value = 27000
$ws.Range("L99").Value = 49999
$ws.Range("M99").Value = -24005
$ws.Range("N99").Value = -55989
$ws.Range("H122").Value = 9056
$ws.Range("J122").Value = 8268.333000000001
$ws.Range("L122").Value = 24804.999
$ws.Range("N122").Value = -29704.999
$ws.Range("H132").Value = 6404.604
$ws.Range("I132").Value = 3843.25
$ws.Range("J132").Value = 9990.5
$ws.Range("K132").Value = 11529.75
$ws.Range("L132").Value = 29971.5
$ws.Range("M132").Value = -8999.75
$ws.Range("N132").Value = -35031.5
$ws.Range("H136").Value = 5832.7617
$ws.Range("I136").Value = 4226.357
$ws.Range("J136").Value = 9045.571
$ws.Range("K136").Value = 12679.071
$ws.Range("L136").Value = 27136.713
$ws.Range("M136").Value = -10129.071
$ws.Range("N136").Value = -32236.713

# Sheet: WVR (30 cell updates)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H53").Value = 48033
$ws.Range("J53").Value = 48033
$ws.Range("L53").Value = 48033
$ws.Range("N53").Value = -49247
$ws.Range("H100").Value = 1123.3438
$ws.Range("I100").Value = 1125.862
$ws.Range("K100").Value = 2251.724
$ws.Range("M100").Value = -1710.724
$ws.Range("H107").Value = 1694.6666
$ws.Range("I107").Value = 1743.1666
$ws.Range("J107").Value = 1597.6666
$ws.Range("K107").Value = 5229.4998
$ws.Range("L107").Value = 4792.9998
$ws.Range("M107").Value = -3309.4998
$ws.Range("N107").Value = -8632.9998
$ws.Range("H122").Value = 1760.25
$ws.Range("I122").Value = 1016.4
$ws.Range("K122").Value = 3049.2
$ws.Range("M122").Value = -599.1999999999998
$ws.Range("H132").Value = 25643256
$ws.Range("I132").Value = 71430344
$ws.Range("J132").Value = 2488.44
$ws.Range("K132").Value = 214291032
$ws.Range("L132").Value = 7465.32
$ws.Range("M132").Value = -214288502
$ws.Range("N132").Value = -12525.32
$ws.Range("H136").Value = 5633.551
$ws.Range("J136").Value = 5623.6333
$ws.Range("L136").Value = 16870.8999
$ws.Range("N136").Value = -21970.8999
